$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the Attribute column (C) for the NVIS rows from "NA" to "object_label"
$ws.Range("C2").Value = "object_label"
$ws.Range("C3").Value = "object_label"
$ws.Range("C4").Value = "object_label"

# Rename the dataset in A4 from "Terrestrial_NEAP" to "NVIS_NEAP"
$ws.Range("A4").Value = "NVIS_NEAP"

# Move the saved selection to B14, matching the author's last cursor position
$ws.Range("B14").Select()
